$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 604
$ws.Range("F4").Value = 893
$ws.Range("F5").Value = 646
$ws.Range("F6").Value = 796
$ws.Range("F7").Value = 367
$ws.Range("F8").Value = 567
$ws.Range("F10").Value = 1141
$ws.Range("F11").Value = 585
$ws.Range("F12").Value = 347
$ws.Range("F13").Value = 464
$ws.Range("F15").Value = 301
$ws.Range("F17").Value = 67
$ws.Range("F18").Value = 529
$ws.Range("F19").Value = 29
$ws.Range("F20").Value = 531
$ws.Range("F22").Value = 503

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 306
$ws.Range("F9").Value = 198

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 306
$ws.Range("F7").Value = 604
$ws.Range("F8").Value = 893
$ws.Range("F9").Value = 646
$ws.Range("F10").Value = 796
$ws.Range("F11").Value = 367
$ws.Range("F12").Value = 567
$ws.Range("F14").Value = 1141
$ws.Range("F15").Value = 585
$ws.Range("F18").Value = 347
$ws.Range("F19").Value = 464
$ws.Range("F23").Value = 301
$ws.Range("F25").Value = 67
$ws.Range("F26").Value = 198
$ws.Range("F28").Value = 529
$ws.Range("F31").Value = 29
$ws.Range("F32").Value = 531
$ws.Range("F34").Value = 503
